$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# The "Motor" quantities in rows 3 and 5 both go from 2 to 4.
# F3 (=D3*E3), F5 (shared formula D5*E5), F6 (=SUM(F3:F5)) and
# F31 (=F6+F13+F24) all recalculate automatically.
$ws.Range("E3").Value = 4
$ws.Range("E5").Value = 4

# Update the window view/selection to match the saved workbook state:
# scrolled so row 4 is at the top, with K18 as the active selected cell.
$ws.Range("K18").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
